$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.342.67'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '1.883.78'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  -0.74%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '245.42'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -2.99%  '
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '43.63'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +5.69%  '
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '53.48'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.0741'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.0971'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -1.18%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '13.40'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +3.24%  '
$ws.Range("D14").Value = '2.155.85'
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.761'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +3.95%  '
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").Value = '1.897.50'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '35.476.93'
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '73.01'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '0.0₃0822'
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '244.58'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '12.81'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '4.94'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -1.62%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.65'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +9.09%  '
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.15'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -6.09%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '165.56'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '8.52'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '18.31'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("D31").Value = '4.128.46'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.71'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +8.49%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '4.27'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.0588'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -3.08%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.90'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -7.02%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '4.15'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.847'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -0.62%  '
$ws.Range("E39").Value = '  -2.17%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.0698'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +7.50%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '17.33'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.0218'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '96.17'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -6.17%  '
$ws.Range("E44").Value = '  -1.85%  '
$ws.Range("D45").Value = '1.306.43'
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '2.33'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -3.47%  '
$ws.Range("E47").Value = '  +7.03%  '
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '6.24'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -5.18%  '
